# Update AVTA course Excel workbook: add new course rows (civil construction
# design + IT telecommunications courses and their package) to the "courses"
# sheet, and refresh the sheet selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("courses")

# --- Row 2: DIPLOMA OF CIVIL CONSTRUCTION DESIGN ---------------------------
$ws.Range("A2").Value = "RII50520"
$ws.Range("B2").Value = "111827M"
$ws.Range("C2").Value = "CIVIL CONSTRUCTION DESIGN"
$ws.Range("D2").Value = "DIPLOMA OF CIVIL CONSTRUCTION DESIGN"
$ws.Range("E2").Value = 52
$ws.Range("H2").Value = "44 wks Tuition + 8 wks Break"
$ws.Range("H2").WrapText = $true
$ws.Range("I2").Value = 19200
$ws.Range("I2").NumberFormat = "#,##0"
$ws.Range("J2").Value = "19,000 tuition fee + 200 handling fee"
$ws.Range("J2").NumberFormat = "#,##0"
$ws.Range("J2").WrapText = $true
$ws.Range("M2").Value = "TAS"

# --- Row 3: ADVANCED DIPLOMA OF CIVIL CONSTRUCTION DESIGN ------------------
$ws.Range("A3").Value = "RII60520"
$ws.Range("B3").Value = "111826A"
$ws.Range("C3").Value = "CIVIL CONSTRUCTION DESIGN"
$ws.Range("D3").Value = "ADVANCED DIPLOMA OF CIVIL CONSTRUCTION DESIGN"
$ws.Range("E3").Value = 104
$ws.Range("H3").Value = "88 wks tuition + 16 wks break"
$ws.Range("H3").WrapText = $true
$ws.Range("I3").Value = 29200
$ws.Range("I3").NumberFormat = "#,##0"
$ws.Range("J3").Value = "29,000 tuition fee + 200 handling fee"
$ws.Range("J3").NumberFormat = "#,##0"
$ws.Range("J3").WrapText = $true
$ws.Range("M3").Value = "TAS"

# --- Row 4: ADVANCED DIPLOMA OF INFORMATION TECHNOLOGY (TELECOMMS) ---------
$ws.Range("A4").Value = "ICT60220"
$ws.Range("B4").Value = "111825B"
$ws.Range("C4").Value = "INFORMATION TECHNOLOGY"
$ws.Range("D4").Value = "ADVANCED DIPLOMA OF INFORMATION TECHNOLOGY `n(TELECOMMUNICATIONS NETWORK ENGINEERING) "
$ws.Range("D4").WrapText = $true
$ws.Range("E4").Value = 104
$ws.Range("H4").Value = "88 wks tuition + 16 wks break"
$ws.Range("H4").WrapText = $true
$ws.Range("I4").Value = 14200
$ws.Range("I4").NumberFormat = "#,##0"
$ws.Range("J4").Value = "14,000 tuition fee + 200 handling fee"
$ws.Range("J4").NumberFormat = "#,##0"
$ws.Range("J4").WrapText = $true
$ws.Range("M4").Value = "TAS"

# --- Row 5: PACKAGES (diploma + advanced diploma civil construction) -------
$ws.Range("A5").Value = "RII50520/RII60520"
$ws.Range("A5").WrapText = $true
$ws.Range("B5").Value = "111827M/111826A"
$ws.Range("B5").WrapText = $true
$ws.Range("C5").Value = "PACKAGES"
$ws.Range("D5").Value = "DIPLOMA OF CIVIL CONSTRUCTION DESIGN + ADVANCED DIPLOMA OF CIVIL CONSTRUCTION DESIGN"
$ws.Range("D5").WrapText = $true
$ws.Range("E5").Value = 104
$ws.Range("H5").Value = "88 wks tuition + 16 wks break"
$ws.Range("H5").WrapText = $true
$ws.Range("I5").Value = 29200
$ws.Range("I5").NumberFormat = "#,##0"
$ws.Range("J5").Value = "29,000 tuition fee + 200 handling fee"
$ws.Range("J5").NumberFormat = "#,##0"
$ws.Range("J5").WrapText = $true
$ws.Range("M5").Value = "TAS"

# Every new data row renders at 45pt tall (wrapped multi-line content).
$ws.Range("A2:R5").RowHeight = 45

# Leave the cursor where the author left it after entering the new rows.
$ws.Range("I18").Select()
